$wb = $excel.ActiveWorkbook

# ---- Sheet: ALERTS ----
$ws = $wb.Worksheets.Item("ALERTS")
$ws.Cells.Item(13,1).Value = "'2026-01-30"
$ws.Cells.Item(13,1).Style = "Normal"
$ws.Cells.Item(13,2).Value = "15:24:41"
$ws.Cells.Item(13,3).Value = "15:00"
$ws.Cells.Item(13,4).Value = "Living Room"
$ws.Cells.Item(13,5).Value = "CRITICAL EMERGENCY"
$ws.Cells.Item(13,6).Value = "FALL_DETECTED"
$ws.Cells.Item(14,1).Value = "'2026-01-30"
$ws.Cells.Item(14,1).Style = "Normal"
$ws.Cells.Item(14,2).Value = "15:24:44"
$ws.Cells.Item(14,3).Value = "15:00"
$ws.Cells.Item(14,4).Value = "Living Room"
$ws.Cells.Item(14,5).Value = "CRITICAL EMERGENCY"
$ws.Cells.Item(14,6).Value = "FALL_DETECTED"
$ws.Cells.Item(15,1).Value = "'2026-01-30"
$ws.Cells.Item(15,1).Style = "Normal"
$ws.Cells.Item(15,2).Value = "15:25:06"
$ws.Cells.Item(15,3).Value = "15:00"
$ws.Cells.Item(15,4).Value = "Living Room"
$ws.Cells.Item(15,5).Value = "CRITICAL EMERGENCY"
$ws.Cells.Item(15,6).Value = "FALL_DETECTED"
$ws.Cells.Item(16,1).Value = "'2026-01-30"
$ws.Cells.Item(16,1).Style = "Normal"
$ws.Cells.Item(16,2).Value = "15:27:01"
$ws.Cells.Item(16,3).Value = "15:00"
$ws.Cells.Item(16,4).Value = "Living Room"
$ws.Cells.Item(16,5).Value = "CRITICAL EMERGENCY"
$ws.Cells.Item(16,6).Value = "FALL_DETECTED"

# ---- Sheet: mmWave ----
$ws = $wb.Worksheets.Item("mmWave")
$ws.Cells.Item(87,1).Value = "'2026-01-30"
$ws.Cells.Item(87,1).Style = "Normal"
$ws.Cells.Item(87,2).Value = "15:25:23"
$ws.Cells.Item(87,3).Value = "15:00"
$ws.Cells.Item(87,4).Value = "Living Room"
$ws.Cells.Item(87,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(87,6).Value = "Active"
$ws.Cells.Item(88,1).Value = "'2026-01-30"
$ws.Cells.Item(88,1).Style = "Normal"
$ws.Cells.Item(88,2).Value = "15:25:33"
$ws.Cells.Item(88,3).Value = "15:00"
$ws.Cells.Item(88,4).Value = "Living Room"
$ws.Cells.Item(88,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(88,6).Value = "Active"
$ws.Cells.Item(89,1).Value = "'2026-01-30"
$ws.Cells.Item(89,1).Style = "Normal"
$ws.Cells.Item(89,2).Value = "15:25:43"
$ws.Cells.Item(89,3).Value = "15:00"
$ws.Cells.Item(89,4).Value = "Living Room"
$ws.Cells.Item(89,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(89,6).Value = "Active"
$ws.Cells.Item(90,1).Value = "'2026-01-30"
$ws.Cells.Item(90,1).Style = "Normal"
$ws.Cells.Item(90,2).Value = "15:25:54"
$ws.Cells.Item(90,3).Value = "15:00"
$ws.Cells.Item(90,4).Value = "Living Room"
$ws.Cells.Item(90,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(90,6).Value = "Active"
$ws.Cells.Item(91,1).Value = "'2026-01-30"
$ws.Cells.Item(91,1).Style = "Normal"
$ws.Cells.Item(91,2).Value = "15:26:04"
$ws.Cells.Item(91,3).Value = "15:00"
$ws.Cells.Item(91,4).Value = "Living Room"
$ws.Cells.Item(91,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(91,6).Value = "Active"
$ws.Cells.Item(92,1).Value = "'2026-01-30"
$ws.Cells.Item(92,1).Style = "Normal"
$ws.Cells.Item(92,2).Value = "15:26:15"
$ws.Cells.Item(92,3).Value = "15:00"
$ws.Cells.Item(92,4).Value = "Living Room"
$ws.Cells.Item(92,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(92,6).Value = "Active"
$ws.Cells.Item(93,1).Value = "'2026-01-30"
$ws.Cells.Item(93,1).Style = "Normal"
$ws.Cells.Item(93,2).Value = "15:26:25"
$ws.Cells.Item(93,3).Value = "15:00"
$ws.Cells.Item(93,4).Value = "Living Room"
$ws.Cells.Item(93,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(93,6).Value = "Active"
$ws.Cells.Item(94,1).Value = "'2026-01-30"
$ws.Cells.Item(94,1).Style = "Normal"
$ws.Cells.Item(94,2).Value = "15:26:36"
$ws.Cells.Item(94,3).Value = "15:00"
$ws.Cells.Item(94,4).Value = "Living Room"
$ws.Cells.Item(94,5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(94,6).Value = "Active"
$ws.Cells.Item(95,1).Value = "'2026-01-30"
$ws.Cells.Item(95,1).Style = "Normal"
$ws.Cells.Item(95,2).Value = "15:30:52"
$ws.Cells.Item(95,3).Value = "15:00"
$ws.Cells.Item(95,4).Value = "Living Room"
$ws.Cells.Item(95,5).Value = "FALL_DETECTED"
$ws.Cells.Item(95,6).Value = "EMERGENCY"
$ws.Cells.Item(96,1).Value = "'2026-01-30"
$ws.Cells.Item(96,1).Style = "Normal"
$ws.Cells.Item(96,2).Value = "15:30:53"
$ws.Cells.Item(96,3).Value = "15:00"
$ws.Cells.Item(96,4).Value = "Living Room"
$ws.Cells.Item(96,5).Value = "FALL_DETECTED"
$ws.Cells.Item(96,6).Value = "EMERGENCY"

# ---- Sheet: Proximity ----
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(38,1).Value = "'2026-01-30"
$ws.Cells.Item(38,1).Style = "Normal"
$ws.Cells.Item(38,2).Value = "15:30:55"
$ws.Cells.Item(38,3).Value = "15:00"
$ws.Cells.Item(38,4).Value = "Living Room Main Door"
$ws.Cells.Item(38,5).Value = "ENTER"
$ws.Cells.Item(38,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(39,1).Value = "'2026-01-30"
$ws.Cells.Item(39,1).Style = "Normal"
$ws.Cells.Item(39,2).Value = "15:30:57"
$ws.Cells.Item(39,3).Value = "15:00"
$ws.Cells.Item(39,4).Value = "Living Room Main Door"
$ws.Cells.Item(39,5).Value = "EXIT"
$ws.Cells.Item(39,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(40,1).Value = "'2026-01-30"
$ws.Cells.Item(40,1).Style = "Normal"
$ws.Cells.Item(40,2).Value = "15:31:00"
$ws.Cells.Item(40,3).Value = "15:00"
$ws.Cells.Item(40,4).Value = "Living Room Main Door"
$ws.Cells.Item(40,5).Value = "ENTER"
$ws.Cells.Item(40,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(41,1).Value = "'2026-01-30"
$ws.Cells.Item(41,1).Style = "Normal"
$ws.Cells.Item(41,2).Value = "15:31:05"
$ws.Cells.Item(41,3).Value = "15:00"
$ws.Cells.Item(41,4).Value = "Living Room Main Door"
$ws.Cells.Item(41,5).Value = "EXIT"
$ws.Cells.Item(41,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(42,1).Value = "'2026-01-30"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).Value = "15:31:11"
$ws.Cells.Item(42,3).Value = "15:00"
$ws.Cells.Item(42,4).Value = "Living Room Main Door"
$ws.Cells.Item(42,5).Value = "ENTER"
$ws.Cells.Item(42,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(43,1).Value = "'2026-01-30"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).Value = "15:31:14"
$ws.Cells.Item(43,3).Value = "15:00"
$ws.Cells.Item(43,4).Value = "Living Room Main Door"
$ws.Cells.Item(43,5).Value = "EXIT"
$ws.Cells.Item(43,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(44,1).Value = "'2026-01-30"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).Value = "15:31:18"
$ws.Cells.Item(44,3).Value = "15:00"
$ws.Cells.Item(44,4).Value = "Living Room Main Door"
$ws.Cells.Item(44,5).Value = "ENTER"
$ws.Cells.Item(44,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(45,1).Value = "'2026-01-30"
$ws.Cells.Item(45,1).Style = "Normal"
$ws.Cells.Item(45,2).Value = "15:31:22"
$ws.Cells.Item(45,3).Value = "15:00"
$ws.Cells.Item(45,4).Value = "Living Room Main Door"
$ws.Cells.Item(45,5).Value = "EXIT"
$ws.Cells.Item(45,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(46,1).Value = "'2026-01-30"
$ws.Cells.Item(46,1).Style = "Normal"
$ws.Cells.Item(46,2).Value = "15:31:24"
$ws.Cells.Item(46,3).Value = "15:00"
$ws.Cells.Item(46,4).Value = "Living Room Main Door"
$ws.Cells.Item(46,5).Value = "ENTER"
$ws.Cells.Item(46,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(47,1).Value = "'2026-01-30"
$ws.Cells.Item(47,1).Style = "Normal"
$ws.Cells.Item(47,2).Value = "15:31:29"
$ws.Cells.Item(47,3).Value = "15:00"
$ws.Cells.Item(47,4).Value = "Living Room Main Door"
$ws.Cells.Item(47,5).Value = "EXIT"
$ws.Cells.Item(47,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(48,1).Value = "'2026-01-30"
$ws.Cells.Item(48,1).Style = "Normal"
$ws.Cells.Item(48,2).Value = "15:31:36"
$ws.Cells.Item(48,3).Value = "15:00"
$ws.Cells.Item(48,4).Value = "Living Room Main Door"
$ws.Cells.Item(48,5).Value = "ENTER"
$ws.Cells.Item(48,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(49,1).Value = "'2026-01-30"
$ws.Cells.Item(49,1).Style = "Normal"
$ws.Cells.Item(49,2).Value = "15:31:40"
$ws.Cells.Item(49,3).Value = "15:00"
$ws.Cells.Item(49,4).Value = "Living Room Main Door"
$ws.Cells.Item(49,5).Value = "EXIT"
$ws.Cells.Item(49,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(50,1).Value = "'2026-01-30"
$ws.Cells.Item(50,1).Style = "Normal"
$ws.Cells.Item(50,2).Value = "15:31:43"
$ws.Cells.Item(50,3).Value = "15:00"
$ws.Cells.Item(50,4).Value = "Living Room Main Door"
$ws.Cells.Item(50,5).Value = "ENTER"
$ws.Cells.Item(50,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(51,1).Value = "'2026-01-30"
$ws.Cells.Item(51,1).Style = "Normal"
$ws.Cells.Item(51,2).Value = "15:31:47"
$ws.Cells.Item(51,3).Value = "15:00"
$ws.Cells.Item(51,4).Value = "Living Room Main Door"
$ws.Cells.Item(51,5).Value = "EXIT"
$ws.Cells.Item(51,6).Value = "User EXITED Living Room Main Door"

# ---- Sheet: Camera ----
$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(37,1).Value = "'2026-01-30"
$ws.Cells.Item(37,1).Style = "Normal"
$ws.Cells.Item(37,2).Value = "15:30:55"
$ws.Cells.Item(37,3).Value = "15:00"
$ws.Cells.Item(37,4).Value = "Living Room Main Door"
$ws.Cells.Item(37,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(37,6).Value = "Active"
$ws.Cells.Item(38,1).Value = "'2026-01-30"
$ws.Cells.Item(38,1).Style = "Normal"
$ws.Cells.Item(38,2).Value = "15:30:57"
$ws.Cells.Item(38,3).Value = "15:00"
$ws.Cells.Item(38,4).Value = "Living Room Main Door"
$ws.Cells.Item(38,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(38,6).Value = "Active"
$ws.Cells.Item(39,1).Value = "'2026-01-30"
$ws.Cells.Item(39,1).Style = "Normal"
$ws.Cells.Item(39,2).Value = "15:31:00"
$ws.Cells.Item(39,3).Value = "15:00"
$ws.Cells.Item(39,4).Value = "Living Room Main Door"
$ws.Cells.Item(39,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(39,6).Value = "Active"
$ws.Cells.Item(40,1).Value = "'2026-01-30"
$ws.Cells.Item(40,1).Style = "Normal"
$ws.Cells.Item(40,2).Value = "15:31:05"
$ws.Cells.Item(40,3).Value = "15:00"
$ws.Cells.Item(40,4).Value = "Living Room Main Door"
$ws.Cells.Item(40,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(40,6).Value = "Active"
$ws.Cells.Item(41,1).Value = "'2026-01-30"
$ws.Cells.Item(41,1).Style = "Normal"
$ws.Cells.Item(41,2).Value = "15:31:10"
$ws.Cells.Item(41,3).Value = "15:00"
$ws.Cells.Item(41,4).Value = "Living Room Main Door"
$ws.Cells.Item(41,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(41,6).Value = "Active"
$ws.Cells.Item(42,1).Value = "'2026-01-30"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).Value = "15:31:14"
$ws.Cells.Item(42,3).Value = "15:00"
$ws.Cells.Item(42,4).Value = "Living Room Main Door"
$ws.Cells.Item(42,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(42,6).Value = "Active"
$ws.Cells.Item(43,1).Value = "'2026-01-30"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).Value = "15:31:18"
$ws.Cells.Item(43,3).Value = "15:00"
$ws.Cells.Item(43,4).Value = "Living Room Main Door"
$ws.Cells.Item(43,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(43,6).Value = "Active"
$ws.Cells.Item(44,1).Value = "'2026-01-30"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).Value = "15:31:22"
$ws.Cells.Item(44,3).Value = "15:00"
$ws.Cells.Item(44,4).Value = "Living Room Main Door"
$ws.Cells.Item(44,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(44,6).Value = "Active"
$ws.Cells.Item(45,1).Value = "'2026-01-30"
$ws.Cells.Item(45,1).Style = "Normal"
$ws.Cells.Item(45,2).Value = "15:31:24"
$ws.Cells.Item(45,3).Value = "15:00"
$ws.Cells.Item(45,4).Value = "Living Room Main Door"
$ws.Cells.Item(45,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(45,6).Value = "Active"
$ws.Cells.Item(46,1).Value = "'2026-01-30"
$ws.Cells.Item(46,1).Style = "Normal"
$ws.Cells.Item(46,2).Value = "15:31:29"
$ws.Cells.Item(46,3).Value = "15:00"
$ws.Cells.Item(46,4).Value = "Living Room Main Door"
$ws.Cells.Item(46,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(46,6).Value = "Active"
$ws.Cells.Item(47,1).Value = "'2026-01-30"
$ws.Cells.Item(47,1).Style = "Normal"
$ws.Cells.Item(47,2).Value = "15:31:36"
$ws.Cells.Item(47,3).Value = "15:00"
$ws.Cells.Item(47,4).Value = "Living Room Main Door"
$ws.Cells.Item(47,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(47,6).Value = "Active"
$ws.Cells.Item(48,1).Value = "'2026-01-30"
$ws.Cells.Item(48,1).Style = "Normal"
$ws.Cells.Item(48,2).Value = "15:31:40"
$ws.Cells.Item(48,3).Value = "15:00"
$ws.Cells.Item(48,4).Value = "Living Room Main Door"
$ws.Cells.Item(48,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(48,6).Value = "Active"
$ws.Cells.Item(49,1).Value = "'2026-01-30"
$ws.Cells.Item(49,1).Style = "Normal"
$ws.Cells.Item(49,2).Value = "15:31:43"
$ws.Cells.Item(49,3).Value = "15:00"
$ws.Cells.Item(49,4).Value = "Living Room Main Door"
$ws.Cells.Item(49,5).Value = "Image Captured (ENTER)"
$ws.Cells.Item(49,6).Value = "Active"
$ws.Cells.Item(50,1).Value = "'2026-01-30"
$ws.Cells.Item(50,1).Style = "Normal"
$ws.Cells.Item(50,2).Value = "15:31:47"
$ws.Cells.Item(50,3).Value = "15:00"
$ws.Cells.Item(50,4).Value = "Living Room Main Door"
$ws.Cells.Item(50,5).Value = "Image Captured (EXIT)"
$ws.Cells.Item(50,6).Value = "Active"

# ---- Sheet: PIR ----
$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(26,1).Value = "'2026-01-30"
$ws.Cells.Item(26,1).Style = "Normal"
$ws.Cells.Item(26,2).Value = "15:25:22"
$ws.Cells.Item(26,3).Value = "15:00"
$ws.Cells.Item(26,4).Value = "Living Room"
$ws.Cells.Item(26,5).Value = "RECOVERY_DETECTION"
$ws.Cells.Item(26,6).Value = "Inactive"

